$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.24"
$ws.Range("E2").Value = "'1.39%"
$ws.Range("D3").Value = "'31.52"
$ws.Range("E3").Value = "'0.05%"
$ws.Range("D4").Value = "'5.131"
$ws.Range("E4").Value = "'0.14%"
$ws.Range("D5").Value = "'0.07976"
$ws.Range("E5").Value = "'6.80%"
$ws.Range("D6").Value = "'2.393"
$ws.Range("E6").Value = "'42.26%"
$ws.Range("D7").Value = "'7.937"
$ws.Range("E7").Value = "'2.67%"
$ws.Range("D8").Value = "'3.854"
$ws.Range("E8").Value = "'1.52%"
$ws.Range("D9").Value = "'0.9145"
$ws.Range("E9").Value = "'-1.59%"
$ws.Range("D10").Value = "'0.1741"
$ws.Range("E10").Value = "'2.95%"
$ws.Range("D11").Value = "'0.07340"
$ws.Range("E11").Value = "'2.42%"
$ws.Range("D12").Value = "'0.08064"
$ws.Range("E12").Value = "'1.11%"
$ws.Range("D13").Value = "'0.03091"
$ws.Range("E13").Value = "'1.98%"
$ws.Range("D14").Value = "'0.09954"
$ws.Range("E14").Value = "'0.49%"
$ws.Range("D15").Value = "'0.001516"
$ws.Range("E15").Value = "'1.18%"
$ws.Range("D16").Value = "'0.006091"
$ws.Range("E16").Value = "'-3.49%"
$ws.Range("D17").Value = "'3.497"
$ws.Range("E17").Value = "'1.17%"
$ws.Range("D18").Value = "'2.239"
$ws.Range("E18").Value = "'0.57%"
$ws.Range("D19").Value = "'0.3253"
$ws.Range("E19").Value = "'-0.88%"
$ws.Range("D20").Value = "'0.1353"
$ws.Range("E20").Value = "'0.28%"
$ws.Range("D21").Value = "'4.686"
$ws.Range("E21").Value = "'2.56%"
$ws.Range("D22").Value = "'0.1608"
$ws.Range("E22").Value = "'3.47%"
$ws.Range("D23").Value = "'0.04627"
$ws.Range("E23").Value = "'-0.37%"
$ws.Range("E24").Value = "'4.12%"
$ws.Range("D25").Value = "'0.004463"
$ws.Range("E25").Value = "'0.98%"
$ws.Range("D26").Value = "'0.0001196"
$ws.Range("E26").Value = "'-8.16%"
$ws.Range("D27").Value = "'0.0003447"
$ws.Range("E27").Value = "'83.43%"
$ws.Range("D39").Value = "'0.01852"
$ws.Range("E39").Value = "'11.56%"
$ws.Range("D40").Value = "'0.04529"
$ws.Range("E40").Value = "'1.98%"
$ws.Range("D41").Value = "'0.007311"
$ws.Range("E41").Value = "'3.71%"
$ws.Range("D42").Value = "'0.1345"
$ws.Range("E42").Value = "'1.51%"
$ws.Range("D43").Value = "'0.002176"
$ws.Range("E43").Value = "'4.95%"
$ws.Range("D44").Value = "'0.01064"
$ws.Range("E44").Value = "'-13.64%"
$ws.Range("D45").Value = "'0.00006406"
$ws.Range("E45").Value = "'6.83%"
$ws.Range("D46").Value = "'0.00000000754"
$ws.Range("E46").Value = "'0.53%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "'0.8206"
$ws.Range("E47").Value = "'-57.22%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.006677"
$ws.Range("E48").Value = "'-39.41%"
$ws.Range("D49").Value = "'0.00002111"
$ws.Range("E49").Value = "'0.53%"
$ws.Range("D50").Value = "'0.0002011"
$ws.Range("E50").Value = "'0.60%"
